$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zero-out the previously populated metadata-compliance / completeness
# counters so the "operation" result starts from a clean slate.
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0

$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0

$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0

$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0

$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0

$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0

# Append a new, empty metadata row (row 14) with the same layout as the
# other rows: attribute name left blank and every numeric counter/score
# column initialised to 0.
$ws.Range("A14").Value = ""
$columns = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $columns) {
    $ws.Range($col + "14").Value = 0
}
